$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.345.65"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "2.511.56"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'176.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "2.511.09"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'4.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "3.023.39"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "'25.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "68.346.39"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "2.531.17"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'11.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'350.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "'4.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").Value = "2.635.12"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0898"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "'508.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'161.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D39").Value = "'18.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'4.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "'150.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.58%  "
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0741"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "'1.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.577"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
